# Updated cryptos list (prices/volumes) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "90.962.73"
$ws.Range("E2").Value = "  +3.35%  "

# Row 3
$ws.Range("D3").Value = "3.070.24"
$ws.Range("E3").Value = "  -1.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.36%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.373"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.60%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.882"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.76%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.13%  "

# Row 10
$ws.Range("D10").Value = "3.068.41"
$ws.Range("E10").Value = "  -1.18%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.681"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +21.68%  "

# Row 12
$ws.Range("E12").Value = "  +5.44%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.91%  "

# Row 14
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "90.743.67"
$ws.Range("E14").Value = "  +2.87%  "

# Row 15
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.48%  "

# Row 17
$ws.Range("D17").Value = "3.630.77"
$ws.Range("E17").Value = "  -1.71%  "

# Row 18
$ws.Range("D18").Value = "3.127.46"
$ws.Range("E18").Value = "  +0.32%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.26%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000222"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
$ws.Range("E24").Value = "  +4.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.60%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "83.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.55%  "

# Row 28
$ws.Range("E28").Value = "  -2.36%  "

# Row 29
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.29%  "

# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.166"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.66%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.53%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "512.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "

# Row 36
$ws.Range("E36").Value = "  +0.52%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.57%  "

# Row 38
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.135"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.40%  "

# Row 41
$ws.Range("E41").Value = "  -0.23%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.138"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.365"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.67%  "

# Row 45
$ws.Range("E45").Value = "  +1.38%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0718"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.08%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.38%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.37%  "

# Row 49
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000264"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.79%  "

# Row 50
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.20%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "165.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.45%  "
